# Updates cryptos list with latest scraped values (GitHub Actions run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.399.35"
$ws.Range("E2").Value = "  +0.05%  "
$ws.Range("D3").Value = "1.847.52"
$ws.Range("E3").Value = "  -0.09%  "
$ws.Range("D4").Value = "'0.9972"
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("D5").Value = "'240.22"
$ws.Range("E5").Value = "  -0.01%  "
$ws.Range("E6").Value = "  -0.34%  "
$ws.Range("D7").Value = "'0.9994"
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").Value = "'0.07501"
$ws.Range("E8").Value = "  -1.68%  "
$ws.Range("D9").Value = "'0.2905"
$ws.Range("E9").Value = "  -0.01%  "
$ws.Range("D10").Value = "'24.46"
$ws.Range("E10").Value = "  -1.05%  "
$ws.Range("D11").Value = "'0.07742"
$ws.Range("E11").Value = "  -0.03%  "
$ws.Range("D12").Value = "1.847.59"
$ws.Range("E12").Value = "  -2.15%  "
$ws.Range("D13").Value = "'5.001"
$ws.Range("E13").Value = "  -0.60%  "
$ws.Range("D14").Value = "'0.6810"
$ws.Range("E14").Value = "  +0.41%  "
$ws.Range("D15").Value = "'0.00001053"
$ws.Range("E15").Value = "  -0.39%  "
$ws.Range("D16").Value = "'82.22"
$ws.Range("E16").Value = "  -1.24%  "
$ws.Range("D17").Value = "2.105.58"
$ws.Range("E17").Value = "  -3.72%  "
$ws.Range("D18").Value = "'6.183"
$ws.Range("E18").Value = "  +0.40%  "
$ws.Range("D19").Value = "29.452.70"
$ws.Range("E19").Value = "  +0.06%  "
$ws.Range("D20").Value = "'229.82"
$ws.Range("E20").Value = "  +0.94%  "
$ws.Range("D21").Value = "'12.35"
$ws.Range("D22").Value = "'0.9989"
$ws.Range("E22").Value = "  -0.05%  "
$ws.Range("D23").Value = "'7.480"
$ws.Range("E23").Value = "  -0.41%  "
$ws.Range("D24").Value = "'0.9992"
$ws.Range("E24").Value = "  -0.11%  "
$ws.Range("D25").Value = "'159.20"
$ws.Range("E25").Value = "  +0.34%  "
$ws.Range("D26").Value = "'0.1375"
$ws.Range("E26").Value = "  -0.73%  "
$ws.Range("D27").Value = "'8.414"
$ws.Range("E27").Value = "  +0.12%  "
$ws.Range("D28").Value = "'17.57"
$ws.Range("E28").Value = "  -0.68%  "
$ws.Range("D29").Value = "'0.06391"
$ws.Range("E29").Value = "  +14.44%  "
$ws.Range("D30").Value = "'1.417"
$ws.Range("E30").Value = "  +3.12%  "
$ws.Range("D31").Value = "'1.475"
$ws.Range("E31").Value = "  +0.85%  "
$ws.Range("D32").Value = "'4.097"
$ws.Range("E32").Value = "  -0.30%  "
$ws.Range("D33").Value = "'4.102"
$ws.Range("E33").Value = "  +0.78%  "
$ws.Range("D34").Value = "'1.830"
$ws.Range("E34").Value = "  -0.27%  "
$ws.Range("D35").Value = "'1.143"
$ws.Range("E35").Value = "  -1.72%  "
$ws.Range("D36").Value = "'0.6979"
$ws.Range("E36").Value = "  -0.12%  "
$ws.Range("D37").Value = "'2.581"
$ws.Range("E37").Value = "  +0.05%  "
$ws.Range("D38").Value = "1.266.12"
$ws.Range("E38").Value = "  +2.79%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "'0.01833"
$ws.Range("E39").Value = "  +1.56%  "
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").Value = "'2.815"
$ws.Range("E40").Value = "  +3.77%  "
$ws.Range("D41").Value = "'6.680"
$ws.Range("E41").Value = "  +4.80%  "
$ws.Range("D42").Value = "'0.9088"
$ws.Range("E42").Value = "  +0.51%  "
$ws.Range("D43").Value = "'0.9994"
$ws.Range("E43").Value = "  -0.10%  "
$ws.Range("D44").Value = "2.010.12"
$ws.Range("E44").Value = "  -18.31%  "
$ws.Range("D45").Value = "'101.38"
$ws.Range("E45").Value = "  -0.15%  "
$ws.Range("D46").Value = "'66.38"
$ws.Range("E46").Value = "  +0.50%  "
$ws.Range("D47").Value = "'1.736"
$ws.Range("E47").Value = "  +3.43%  "
$ws.Range("D48").Value = "'7.082"
$ws.Range("E48").Value = "  -1.57%  "
$ws.Range("E49").Value = "  +3.41%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'9.071"
$ws.Range("E50").Value = "  +0.53%  "
$ws.Range("B51").Value = "BabyDogeCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D51").Value = "'0.00000000116"
$ws.Range("E51").Value = "  -4.57%  "
